$p = $ppt.ActivePresentation

# --- 1. Clear the "Name: Dhananjay Adik" subtitle text on the title slide ---
$titleSlide = $p.Slides.Item(1)
$subtitle = $titleSlide.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Text = ""

# --- 2. Re-cache the "datetimeFigureOut" date placeholder text (5/21/2021 -> 12/27/2018) ---
# across the slide master and every slide layout.
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "12/27/2018"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}
